$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.550.67"
$ws.Range("E2").Value = "  +2.74%  "

$ws.Range("D3").Value = "1.668.42"
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'237.28"
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("D7").Value = "'0.4738"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D9").Value = "'0.06170"
$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.667.98"
$ws.Range("E10").Value = "  +2.04%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07009"
$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").Value = "'14.77"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("E13").Value = "  -3.31%  "

$ws.Range("D14").Value = "'4.359"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").Value = "'75.30"
$ws.Range("E15").Value = "  +3.44%  "

$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").Value = "25.547.91"
$ws.Range("E18").Value = "  +2.67%  "

$ws.Range("D19").Value = "'0.000006730"
$ws.Range("E19").Value = "  +2.61%  "

$ws.Range("D20").Value = "'11.40"
$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("D21").Value = "1.882.58"
$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("E22").Value = "  +1.98%  "

$ws.Range("D23").Value = "'8.775"
$ws.Range("E23").Value = "  +2.62%  "

$ws.Range("D24").Value = "'5.226"
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").Value = "'137.09"
$ws.Range("E25").Value = "  +2.90%  "

$ws.Range("D26").Value = "'14.98"
$ws.Range("E26").Value = "  +1.61%  "

$ws.Range("D27").Value = "'1.387"
$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("E28").Value = "  +5.27%  "

$ws.Range("D29").Value = "'104.34"
$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("D30").Value = "'3.990"
$ws.Range("E30").Value = "  +6.02%  "

$ws.Range("D31").Value = "'0.07824"

$ws.Range("D32").Value = "'3.620"
$ws.Range("E32").Value = "  +2.67%  "

$ws.Range("B33").Value = "Frax"
$ws.Range("C33").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D33").Value = "'0.9992"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04310"
$ws.Range("E34").Value = "  +0.59%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.624"
$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9524"
$ws.Range("E36").Value = "  +3.64%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6053"
$ws.Range("E37").Value = "  +4.36%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.9411"
$ws.Range("E38").Value = "  +15.26%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.499"
$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "'0.9999"
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").Value = "'1.849"
$ws.Range("E41").Value = "  +4.30%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01475"
$ws.Range("E42").Value = "  -4.02%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.07"
$ws.Range("E43").Value = "  +2.19%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3742"
$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'4.891"
$ws.Range("E45").Value = "  +4.06%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1113"
$ws.Range("E46").Value = "  +2.30%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'6.185"
$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05262"
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'29.84"
$ws.Range("E49").Value = "  +1.42%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.488"
$ws.Range("E50").Value = "  +4.16%  "

$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  +0.14%  "
